$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" (Total) sheet: the existing rows 2/3 get overwritten with
#    the new "2022-Q3"/"2022-Q2" figures, and a brand-new row 4 is
#    appended for "2022-Q1" (mirrors the original file's row layout,
#    where the index column A is just a running 0/1/2 counter).
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Row 3 used to hold "2022-Q1" (16, 7.89); reuse it for "2022-Q2"
# (15, 5.17) -- i.e. what row 2 used to contain.
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 15
$totalSheet.Range("D3").Value = 5.17

# Append a new row 4 (copy row 3's formatting) for "2022-Q1".
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A4:D4").PasteSpecial(-4122)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 16
$totalSheet.Range("D4").Value = 7.89

# Row 2 used to hold "2022-Q2" (15, 5.17); overwrite with the new
# "2022-Q3" figures.
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 18
$totalSheet.Range("D2").Value = 5.89

# ------------------------------------------------------------------
# 2) Add a new "2022-Q3" detail worksheet positioned right before the
#    existing "2022-Q2" sheet (so tab order becomes
#    总计, 2022-Q3, 2022-Q2, 2022-Q1, matching the updated "总计" list).
#    Duplicating "2022-Q2" keeps headers/column formatting identical.
# ------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# The fund-code / decimal-looking text columns must stay text (leading
# zeros like "040047" and trailing zeros like "5.17" must survive) --
# force a text number format before writing any values into them.
$newSheet.Range("B2:B16").NumberFormat = "@"
$newSheet.Range("D2:G16").NumberFormat = "@"

# Extend formatting down to rows 17-19 (the old sheet only had data
# through row 16) by copying row 16's formatting.
$newSheet.Range("A16:H16").Copy()
$newSheet.Range("A17:H19").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "159941"
$newSheet.Range("C2").Value = "广发纳斯达克100ETF（QDII）"
$newSheet.Range("D2").Value = "106.15"
$newSheet.Range("E2").Value = "91.14"
$newSheet.Range("F2").Value = "1.88"
$newSheet.Range("G2").Value = "1.9956"
$newSheet.Range("H2").Value = 10
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "513100"
$newSheet.Range("C3").Value = "国泰纳斯达克100（QDII-ETF）"
$newSheet.Range("D3").Value = "46.54"
$newSheet.Range("E3").Value = "91.35"
$newSheet.Range("F3").Value = "1.88"
$newSheet.Range("G3").Value = "0.8750"
$newSheet.Range("H3").Value = 10
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "040047"
$newSheet.Range("C4").Value = "华安纳斯达克100指数（QDII）美元现钞A"
$newSheet.Range("D4").Value = "24.52"
$newSheet.Range("E4").Value = "92.09"
$newSheet.Range("F4").Value = "1.92"
$newSheet.Range("G4").Value = "0.4708"
$newSheet.Range("H4").Value = 10
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "040048"
$newSheet.Range("C5").Value = "华安纳斯达克100指数（QDII）美元现汇A"
$newSheet.Range("D5").Value = "24.52"
$newSheet.Range("E5").Value = "92.09"
$newSheet.Range("F5").Value = "1.92"
$newSheet.Range("G5").Value = "0.4708"
$newSheet.Range("H5").Value = 10
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "040046"
$newSheet.Range("C6").Value = "华安纳斯达克100指数（QDII）人民币A"
$newSheet.Range("D6").Value = "22.21"
$newSheet.Range("E6").Value = "92.09"
$newSheet.Range("F6").Value = "1.92"
$newSheet.Range("G6").Value = "0.4264"
$newSheet.Range("H6").Value = 10
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "160213"
$newSheet.Range("C7").Value = "国泰纳斯达克100指数（QDII）"
$newSheet.Range("D7").Value = "15.14"
$newSheet.Range("E7").Value = "85.81"
$newSheet.Range("F7").Value = "1.86"
$newSheet.Range("G7").Value = "0.2816"
$newSheet.Range("H7").Value = 10
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "000834"
$newSheet.Range("C8").Value = "大成纳斯达克100指数（QDII）"
$newSheet.Range("D8").Value = "14.15"
$newSheet.Range("E8").Value = "85.22"
$newSheet.Range("F8").Value = "1.77"
$newSheet.Range("G8").Value = "0.2505"
$newSheet.Range("H8").Value = 10
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "000043"
$newSheet.Range("C9").Value = "嘉实美国成长股票（QDII）人民币"
$newSheet.Range("D9").Value = "12.41"
$newSheet.Range("E9").Value = "92.80"
$newSheet.Range("F9").Value = "1.88"
$newSheet.Range("G9").Value = "0.2333"
$newSheet.Range("H9").Value = 6
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "000044"
$newSheet.Range("C10").Value = "嘉实美国成长股票（QDII）美元现汇"
$newSheet.Range("D10").Value = "12.41"
$newSheet.Range("E10").Value = "92.80"
$newSheet.Range("F10").Value = "1.88"
$newSheet.Range("G10").Value = "0.2333"
$newSheet.Range("H10").Value = 6
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "513300"
$newSheet.Range("C11").Value = "华夏纳斯达克100ETF（QDII）"
$newSheet.Range("D11").Value = "11.08"
$newSheet.Range("E11").Value = "97.32"
$newSheet.Range("F11").Value = "2.02"
$newSheet.Range("G11").Value = "0.2238"
$newSheet.Range("H11").Value = 4
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "003722"
$newSheet.Range("C12").Value = "易方达纳斯达克100指数美元（QDII-LOF）A"
$newSheet.Range("D12").Value = "7.72"
$newSheet.Range("E12").Value = "90.67"
$newSheet.Range("F12").Value = "1.91"
$newSheet.Range("G12").Value = "0.1475"
$newSheet.Range("H12").Value = 10
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "161130"
$newSheet.Range("C13").Value = "易方达纳斯达克100指数人民币（QDII-LOF）"
$newSheet.Range("D13").Value = "7.72"
$newSheet.Range("E13").Value = "90.67"
$newSheet.Range("F13").Value = "1.91"
$newSheet.Range("G13").Value = "0.1475"
$newSheet.Range("H13").Value = 10
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "014978"
$newSheet.Range("C14").Value = "华安纳斯达克100指数（QDII）人民币C"
$newSheet.Range("D14").Value = "2.31"
$newSheet.Range("E14").Value = "92.09"
$newSheet.Range("F14").Value = "1.92"
$newSheet.Range("G14").Value = "0.0444"
$newSheet.Range("H14").Value = 10
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "013329"
$newSheet.Range("C15").Value = "嘉实全球价值股票（QDII）美元现汇"
$newSheet.Range("D15").Value = "1.68"
$newSheet.Range("E15").Value = "90.63"
$newSheet.Range("F15").Value = "1.61"
$newSheet.Range("G15").Value = "0.0270"
$newSheet.Range("H15").Value = 7
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "013328"
$newSheet.Range("C16").Value = "嘉实全球价值股票（QDII）人民币"
$newSheet.Range("D16").Value = "1.68"
$newSheet.Range("E16").Value = "90.63"
$newSheet.Range("F16").Value = "1.61"
$newSheet.Range("G16").Value = "0.0270"
$newSheet.Range("H16").Value = 7
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "159632"
$newSheet.Range("C17").Value = "华安纳斯达克100ETF（QDII）"
$newSheet.Range("D17").Value = "1.51"
$newSheet.Range("E17").Value = "89.05"
$newSheet.Range("F17").Value = "1.70"
$newSheet.Range("G17").Value = "0.0257"
$newSheet.Range("H17").Value = 10
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "012871"
$newSheet.Range("C18").Value = "易方达纳斯达克100指数美元（QDII-LOF）C"
$newSheet.Range("D18").Value = "0.18"
$newSheet.Range("E18").Value = "90.67"
$newSheet.Range("F18").Value = "1.91"
$newSheet.Range("G18").Value = "0.0034"
$newSheet.Range("H18").Value = 10
$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "012870"
$newSheet.Range("C19").Value = "易方达纳斯达克100指数人民币（QDII-LOF）C"
$newSheet.Range("D19").Value = "0.18"
$newSheet.Range("E19").Value = "90.67"
$newSheet.Range("F19").Value = "1.91"
$newSheet.Range("G19").Value = "0.0034"
$newSheet.Range("H19").Value = 10

# Restore the originally active sheet ("2022-Q1") -- copying a sheet
# makes the copy the active tab as a side effect.
$wb.Worksheets.Item("2022-Q1").Activate()

Write-Host "Edit complete"
